$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.267.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6036"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07056"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2797"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.18%  "
$ws.Range("E10").Value = "  -5.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07656"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.836.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.796"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6294"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009778"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "79.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.271.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.843"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "224.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.57%  "
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.013"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "156.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.989"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.90%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1300"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.84%  "
$ws.Range("B28").Value = "Hedera"
$ws.Range("C28").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06585"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.74%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.472"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.449"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.851"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.800"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.90%  "
$ws.Range("E33").Value = "  -2.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.721"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6475"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.547"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.214.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01750"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.93%  "
$ws.Range("E40").Value = "  -5.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8982"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.993.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.95%  "
$ws.Range("E46").Value = "  -2.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.582"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.580"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4555"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05502"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.410"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.53%  "
